$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.54598331451416
$ws.Range("B1").Value = 6.247114658355713
$ws.Range("C1").Value = 5.722426891326904
$ws.Range("D1").Value = 6.552399158477783
$ws.Range("E1").Value = 3.84219765663147
